$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1 (09:00 -> 09:15)
$ws.Range("F1").Value = "Last status check on: 29.01.2022 09:15"

# Row 8 (Benzina Albert Modrice) price refresh:
#   B8: new current price
#   C8: old price (previous current price)
#   D8: delta, now stored as a text string with explicit sign
#   E8: last-checked timestamp, now stored as plain text instead of a date serial
$ws.Range("B8").Value = 36.4
$ws.Range("C8").Value = 36.2

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "+0.2"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2022-01-29 09:15:21"
$ws.Range("E8").Style = "Normal"
